# The shared "Additional Guests" value for Humphrey (row 2) is no longer a
# placeholder count of 1 - it now holds the real decimal figure (1.23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.23

# The author's selection when they last saved the sheet was D3, not H12.
$ws.Range("D3").Select()
